$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.85067314828899
$ws.Range("D2").Value = 11.70583913879577
$ws.Range("E2").Value = 14.27620259265731
$ws.Range("F2").Value = 58.23728371935554
$ws.Range("G2").Value = 3.779119303145857
$ws.Range("K2").Value = 24.52572770006956
$ws.Range("L2").Value = 9.572060380437346
$ws.Range("C3").Value = 13.77958305940087
$ws.Range("D3").Value = 11.63094892508382
$ws.Range("E3").Value = 14.24035418084026
$ws.Range("F3").Value = 57.36330371113632
$ws.Range("G3").Value = 3.784890121893208
$ws.Range("K3").Value = 24.32282442107993
$ws.Range("L3").Value = 9.586209900224526
$ws.Range("C4").Value = 13.73972747118372
$ws.Range("D4").Value = 11.58567330129396
$ws.Range("E4").Value = 14.22205670809225
$ws.Range("F4").Value = 56.83115509786165
$ws.Range("G4").Value = 3.788605825818031
$ws.Range("K4").Value = 24.20754881628714
$ws.Range("L4").Value = 9.596698906697124
$ws.Range("C5").Value = 13.72444831557754
$ws.Range("D5").Value = 11.5674035269186
$ws.Range("E5").Value = 14.21553521735095
$ws.Range("F5").Value = 56.61561534877201
$ws.Range("G5").Value = 3.790163581357405
$ws.Range("K5").Value = 24.16296045073804
$ws.Range("L5").Value = 9.601424923188132
$ws.Range("C6").Value = 13.72196955788665
$ws.Range("D6").Value = 11.56438080746714
$ws.Range("E6").Value = 14.21450880557876
$ws.Range("F6").Value = 56.57991002321746
$ws.Range("G6").Value = 3.790424883949386
$ws.Range("K6").Value = 24.15570204254071
$ws.Range("L6").Value = 9.6022369159367
$ws.Range("C7").Value = 13.73951750489003
$ws.Range("D7").Value = 11.58542617519076
$ws.Range("E7").Value = 14.22196497097555
$ws.Range("F7").Value = 56.82824268132601
$ws.Range("G7").Value = 3.788626657516721
$ws.Range("K7").Value = 24.20693775758014
$ws.Range("L7").Value = 9.596760816339305
$ws.Range("C8").Value = 13.82537925375347
$ws.Range("D8").Value = 11.67987185044963
$ws.Range("E8").Value = 14.26307221847948
$ws.Range("F8").Value = 57.93512007887056
$ws.Range("G8").Value = 3.781073436406262
$ws.Range("K8").Value = 24.45386452858575
$ws.Range("L8").Value = 9.576564573130169
$ws.Range("C9").Value = 14.02342681607906
$ws.Range("D9").Value = 11.87056542456591
$ws.Range("E9").Value = 14.37306271781191
$ws.Range("F9").Value = 60.13284572615826
$ws.Range("G9").Value = 3.767619184080336
$ws.Range("K9").Value = 25.00963187690923
$ws.Range("L9").Value = 9.551305388873722
$ws.Range("C10").Value = 14.18635364799029
$ws.Range("D10").Value = 12.01382645067255
$ws.Range("E10").Value = 14.47161091810936
$ws.Range("F10").Value = 61.75281657428474
$ws.Range("G10").Value = 3.758547479728166
$ws.Range("K10").Value = 25.45816083498346
$ws.Range("L10").Value = 9.541564868760586
$ws.Range("C11").Value = 14.26408703142667
$ws.Range("D11").Value = 12.07963885073523
$ws.Range("E11").Value = 14.52023888905549
$ws.Range("F11").Value = 62.48858000176957
$ws.Range("G11").Value = 3.754593957086378
$ws.Range("K11").Value = 25.67013579783285
$ws.Range("L11").Value = 9.539062260670716
$ws.Range("C12").Value = 14.2940267033328
$ws.Range("D12").Value = 12.10464800408007
$ws.Range("E12").Value = 14.53919291839818
$ws.Range("F12").Value = 62.76683361687009
$ws.Range("G12").Value = 3.753121526026335
$ws.Range("K12").Value = 25.75147512722769
$ws.Range("L12").Value = 9.538392871069796
$ws.Range("C13").Value = 14.2875565279191
$ws.Range("D13").Value = 12.09925802854014
$ws.Range("E13").Value = 14.53508693897667
$ws.Range("F13").Value = 62.70692597213653
$ws.Range("G13").Value = 3.753437546324252
$ws.Range("K13").Value = 25.73391087365728
$ws.Range("L13").Value = 9.538524644004937
$ws.Range("C14").Value = 14.26654018839475
$ws.Range("D14").Value = 12.08169464587243
$ws.Range("E14").Value = 14.52178747267073
$ws.Range("F14").Value = 62.51148049214322
$ws.Range("G14").Value = 3.754472325950581
$ws.Range("K14").Value = 25.67680662968532
$ws.Range("L14").Value = 9.539001606780763
$ws.Range("C15").Value = 14.25373216934634
$ws.Range("D15").Value = 12.07094780137712
$ws.Range("E15").Value = 14.51371124039076
$ws.Range("F15").Value = 62.39171136526662
$ws.Range("G15").Value = 3.755109365854728
$ws.Range("K15").Value = 25.64196566972176
$ws.Range("L15").Value = 9.539330030321581
$ws.Range("C16").Value = 14.18134499505196
$ws.Range("D16").Value = 12.00953796422845
$ws.Range("E16").Value = 14.46850885671782
$ws.Range("F16").Value = 61.70469230814918
$ws.Range("G16").Value = 3.758809318178947
$ws.Range("K16").Value = 25.44446134817058
$ws.Range("L16").Value = 9.541767311170366
$ws.Range("C17").Value = 14.13785289748033
$ws.Range("D17").Value = 11.9720260954897
$ws.Range("E17").Value = 14.44174707061081
$ws.Range("F17").Value = 61.28278785463736
$ws.Range("G17").Value = 3.761123327263604
$ws.Range("K17").Value = 25.32528088841382
$ws.Range("L17").Value = 9.543757146424845
$ws.Range("C18").Value = 14.11317890843372
$ws.Range("D18").Value = 11.9505110253703
$ws.Range("E18").Value = 14.42671231862743
$ws.Range("F18").Value = 61.04002092242391
$ws.Range("G18").Value = 3.762470604180904
$ws.Range("K18").Value = 25.25748297641166
$ws.Range("L18").Value = 9.545083103505531
$ws.Range("C19").Value = 14.10488385852524
$ws.Range("D19").Value = 11.9432369742781
$ws.Range("E19").Value = 14.4216834656308
$ws.Range("F19").Value = 60.95781332064624
$ws.Range("G19").Value = 3.762929578960858
$ws.Range("K19").Value = 25.23465901478771
$ws.Range("L19").Value = 9.545563183511414
$ws.Range("C20").Value = 14.14244747566851
$ws.Range("D20").Value = 11.97601304385462
$ws.Range("E20").Value = 14.44455890918249
$ws.Range("F20").Value = 61.32771194958996
$ws.Range("G20").Value = 3.760875309639312
$ws.Range("K20").Value = 25.33789058817491
$ws.Range("L20").Value = 9.543526537565727
$ws.Range("C21").Value = 14.27269966154724
$ws.Range("D21").Value = 12.08685110502453
$ws.Range("E21").Value = 14.52567925750363
$ws.Range("F21").Value = 62.56889898052759
$ws.Range("G21").Value = 3.754167718060484
$ws.Range("K21").Value = 25.69355110323733
$ws.Range("L21").Value = 9.538853950949974
$ws.Range("C22").Value = 14.36075480556086
$ws.Range("D22").Value = 12.159797299782
$ws.Range("E22").Value = 14.5818381253173
$ws.Range("F22").Value = 63.37787788318785
$ws.Range("G22").Value = 3.749927678163144
$ws.Range("K22").Value = 25.93219210151801
$ws.Range("L22").Value = 9.537422655428095
$ws.Range("C23").Value = 14.31349589858248
$ws.Range("D23").Value = 12.12081983832995
$ws.Range("E23").Value = 14.5515799381144
$ws.Range("F23").Value = 62.94637593852593
$ws.Range("G23").Value = 3.752177589541674
$ws.Range("K23").Value = 25.8042822895528
$ws.Range("L23").Value = 9.538037798991956
$ws.Range("C24").Value = 14.14036923714585
$ws.Range("D24").Value = 11.97421038526536
$ws.Range("E24").Value = 14.44328658333674
$ws.Range("F24").Value = 61.30740240059313
$ws.Range("G24").Value = 3.760987385680506
$ws.Range("K24").Value = 25.33218749318654
$ws.Range("L24").Value = 9.543630229097793
$ws.Range("C25").Value = 13.96673668062704
$ws.Range("D25").Value = 11.81840815216304
$ws.Range("E25").Value = 14.34017615439287
$ws.Range("F25").Value = 59.53656133028802
$ws.Range("G25").Value = 3.771115070461174
$ws.Range("K25").Value = 24.85195885945865
$ws.Range("L25").Value = 9.55659514054239
